# Atualização de bases das ligas, do dia: 11-03-2024 às 22:32
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Range("B130").Value = 7453204
$ws.Range("F130").Value = "Cerro Porteno"
$ws.Range("G130").Value = "Tacuary"
$ws.Range("H130").Value = 1
$ws.Range("J130").Value = "D"
$ws.Range("K130").Value = 1.285
$ws.Range("L130").Value = 5
$ws.Range("M130").Value = 8
$ws.Range("N130").Value = 1.285
$ws.Range("O130").Value = 4.75
$ws.Range("P130").Value = 8
$ws.Range("Q130").Value = -1.5
$ws.Range("R130").Value = 1.9
$ws.Range("S130").Value = 1.9
$ws.Range("T130").Value = 3
$ws.Range("U130").Value = 1.9
$ws.Range("V130").Value = 1.9
$ws.Range("X130").Value = 3.75
$ws.Range("Y130").Value = -1
$ws.Range("AA130").Value = 0.8999999999999999
$ws.Range("AC130").Value = 0.8999999999999999

# Row 131
$ws.Range("B131").Value = 7454842
$ws.Range("F131").Value = "Sportivo Luqueno"
$ws.Range("G131").Value = "Libertad Asuncion"
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = "A"
$ws.Range("K131").Value = 4
$ws.Range("L131").Value = 3.6
$ws.Range("M131").Value = 1.727
$ws.Range("N131").Value = 3.5
$ws.Range("O131").Value = 3.3
$ws.Range("P131").Value = 1.95
$ws.Range("Q131").Value = 0.5
$ws.Range("R131").Value = 1.8
$ws.Range("S131").Value = 2
$ws.Range("T131").Value = 2.5
$ws.Range("U131").Value = 1.975
$ws.Range("V131").Value = 1.825
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = 0.95
$ws.Range("AA131").Value = 1
$ws.Range("AC131").Value = 0.825

# Row 134
$ws.Range("B134").Value = 7493428
$ws.Range("F134").Value = "Guairena FC"
$ws.Range("G134").Value = "Resistencia FC"
$ws.Range("H134").Value = 4
$ws.Range("J134").Value = "H"
$ws.Range("K134").Value = 1.727
$ws.Range("L134").Value = 3.6
$ws.Range("M134").Value = 4.2
$ws.Range("N134").Value = 1.45
$ws.Range("O134").Value = 4.2
$ws.Range("P134").Value = 6
$ws.Range("Q134").Value = -1
$ws.Range("R134").Value = 1.775
$ws.Range("S134").Value = 2.025
$ws.Range("T134").Value = 2.75
$ws.Range("U134").Value = 1.825
$ws.Range("V134").Value = 1.975
$ws.Range("W134").Value = 0.45
$ws.Range("X134").Value = -1
$ws.Range("Z134").Value = 0.7749999999999999
$ws.Range("AA134").Value = -1
$ws.Range("AB134").Value = 0.825
$ws.Range("AC134").Value = -1

# Row 135
$ws.Range("B135").Value = 7493427
$ws.Range("F135").Value = "Tacuary"
$ws.Range("G135").Value = "Sportivo Luqueno"
$ws.Range("H135").Value = 1
$ws.Range("J135").Value = "D"
$ws.Range("K135").Value = 3.4
$ws.Range("L135").Value = 3.3
$ws.Range("M135").Value = 2
$ws.Range("N135").Value = 3.2
$ws.Range("O135").Value = 3.25
$ws.Range("P135").Value = 2.1
$ws.Range("Q135").Value = 0.25
$ws.Range("R135").Value = 2.025
$ws.Range("S135").Value = 1.775
$ws.Range("T135").Value = 2.5
$ws.Range("U135").Value = 1.975
$ws.Range("V135").Value = 1.825
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = 2.25
$ws.Range("Z135").Value = 0.5125
$ws.Range("AA135").Value = -0.5
$ws.Range("AB135").Value = -1
$ws.Range("AC135").Value = 0.825

# Row 144
$ws.Range("B144").Value = 7493311
$ws.Range("F144").Value = "General Caballero JLM"
$ws.Range("G144").Value = "Olimpia Asuncion"
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 3.4
$ws.Range("L144").Value = 3.3
$ws.Range("M144").Value = 2
$ws.Range("N144").Value = 3.2
$ws.Range("O144").Value = 3.25
$ws.Range("P144").Value = 2.1
$ws.Range("Q144").Value = 0.25
$ws.Range("R144").Value = 1.95
$ws.Range("S144").Value = 1.85
$ws.Range("T144").Value = 2.25
$ws.Range("U144").Value = 1.775
$ws.Range("V144").Value = 2.025
$ws.Range("W144").Value = -1
$ws.Range("Y144").Value = 1.1
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.8500000000000001
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 1.025

# Row 145
$ws.Range("B145").Value = 7493312
$ws.Range("F145").Value = "Cerro Porteno"
$ws.Range("G145").Value = "Guarani Asuncion"
$ws.Range("H145").Value = 4
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = "H"
$ws.Range("K145").Value = 1.7
$ws.Range("L145").Value = 3.6
$ws.Range("M145").Value = 4.333
$ws.Range("N145").Value = 1.727
$ws.Range("O145").Value = 3.75
$ws.Range("P145").Value = 4.2
$ws.Range("Q145").Value = -0.5
$ws.Range("R145").Value = 1.8
$ws.Range("S145").Value = 2
$ws.Range("T145").Value = 2.75
$ws.Range("U145").Value = 1.875
$ws.Range("V145").Value = 1.925
$ws.Range("W145").Value = 0.7270000000000001
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 0.8
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 0.875
$ws.Range("AC145").Value = -1

# Row 194
$ws.Range("H194").Value = 0
$ws.Range("I194").Value = 1
$ws.Range("J194").Value = "A"
$ws.Range("R194").Value = 1.975
$ws.Range("S194").Value = 1.825
$ws.Range("U194").Value = 1.825
$ws.Range("V194").Value = 1.975
$ws.Range("W194").Value = -1
$ws.Range("X194").Value = -1
$ws.Range("Y194").Value = 7.5
$ws.Range("Z194").Value = -1
$ws.Range("AA194").Value = 0.825
$ws.Range("AB194").Value = -1
$ws.Range("AC194").Value = 0.9750000000000001

# Row 195
$ws.Range("H195").Value = 1
$ws.Range("I195").Value = 1
$ws.Range("J195").Value = "D"
$ws.Range("N195").Value = 5.25
$ws.Range("O195").Value = 3.5
$ws.Range("P195").Value = 1.615
$ws.Range("Q195").Value = 0.75
$ws.Range("R195").Value = 2
$ws.Range("S195").Value = 1.8
$ws.Range("U195").Value = 1.8
$ws.Range("V195").Value = 2
$ws.Range("W195").Value = -1
$ws.Range("X195").Value = 2.5
$ws.Range("Y195").Value = -1
$ws.Range("Z195").Value = 1
$ws.Range("AA195").Value = -1
$ws.Range("AB195").Value = -1
$ws.Range("AC195").Value = 1

# Row 196
$ws.Range("H196").Value = 1
$ws.Range("I196").Value = 0
$ws.Range("J196").Value = "H"
$ws.Range("W196").Value = 1.1
$ws.Range("X196").Value = -1
$ws.Range("Y196").Value = -1
$ws.Range("Z196").Value = 0.825
$ws.Range("AA196").Value = -1
$ws.Range("AB196").Value = -1
$ws.Range("AC196").Value = 0.825

# Row 197 removed entirely (match no longer present in source feed)
$ws.Rows(197).Delete()
